$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# =====================================================================
# Row 1 : green title banner "Solid driver check result (RPMs)"
# =====================================================================
$ws.Range("A1").Value = "Solid driver check result (RPMs)"
$ws.Range("A1:J1").Merge()

$titleCell = $ws.Range("A1")
$titleCell.Font.Size = 18
$titleCell.Font.Bold = $true
$titleCell.Interior.Color = 7911984
$titleCell.HorizontalAlignment = -4108

# =====================================================================
# Row 2 : blank spacer row (Calibri 11)
# =====================================================================
$spacer2 = $ws.Range("A2:J2")
$spacer2.Font.Name = "Calibri"

# =====================================================================
# Row 3 : description paragraph
# =====================================================================
$ws.Range("A3").Value = "soliddriver-checks is a tool for parnter(s) and customer(s) to check their RPMs to ensure these are meet basic SUSE requirements."
$ws.Range("A3:J3").Merge()
$descCell = $ws.Range("A3")
$descCell.Font.Name = "Poppins"
$descCell.Font.Size = 14
$descCell.WrapText = $true

# =====================================================================
# Row 4 : blank spacer row (Calibri 11)
# =====================================================================
$spacer4 = $ws.Range("A4:J4")
$spacer4.Font.Name = "Calibri"

# =====================================================================
# Row 5 : "Please refer to Kernel Module Packages Manual..." paragraph
# =====================================================================
$ws.Range("A5").Value = "Please refer to Kernel Module Packages Manual to learn how to build a KMP(Kernel Module Package)."
$ws.Range("A5:J5").Merge()
$kmpCell = $ws.Range("A5")
$kmpCell.Font.Name = '"poppins medium"'
$kmpCell.WrapText = $true

# =====================================================================
# Row 6 : blank spacer row (Calibri 11)
# =====================================================================
$spacer6 = $ws.Range("A6:J6")
$spacer6.Font.Name = "Calibri"

# =====================================================================
# Row 7 : "What do we check?" table header (boxed)
# =====================================================================
$ws.Range("A7").Value = "What do we check?"
$headCell = $ws.Range("A7")
$headCell.Font.Name = '"poppins medium"'
$headCell.Borders.Item(7).LineStyle = 1
$headCell.Borders.Item(8).LineStyle = 1
$headCell.Borders.Item(9).LineStyle = 1

$midHead = $ws.Range("B7:I7")
$midHead.Borders.Item(8).LineStyle = 1
$midHead.Borders.Item(9).LineStyle = 1

$rightHead = $ws.Range("J7")
$rightHead.Borders.Item(8).LineStyle = 1
$rightHead.Borders.Item(9).LineStyle = 1
$rightHead.Borders.Item(10).LineStyle = 1

